$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.690.76'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '3.406.74'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E5').Value = '  -0.20%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '654.86'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('E7').Value = '  +1.01%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.436'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +5.54%  '
$ws.Range('E9').Value = '  +5.68%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').Value = '3.401.10'
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('E12').Value = '  +4.03%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '41.74'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.09%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.42'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +19.19%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0000261'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +3.41%  '
$ws.Range('D16').Value = '97.439.16'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '4.039.96'
$ws.Range('E17').Value = '  +2.61%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '8.62'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +32.95%  '
$ws.Range('D19').Value = '3.409.96'
$ws.Range('E19').Value = '  +2.96%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '17.55'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +10.16%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.498'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +43.62%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '3.46'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.55%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.71'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +12.99%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '507.09'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +3.96%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.0000207'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.01%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '6.20'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +6.62%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '99.00'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +10.68%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '12.80'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +4.98%  '
$ws.Range('D29').Value = '3.586.82'
$ws.Range('E29').Value = '  +2.73%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.154'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.35%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.202'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +5.88%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '11.44'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +6.93%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.18%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.08%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.568'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +17.38%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '29.73'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +5.63%  '
$ws.Range('E37').Value = '  +16.09%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '7.73'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +4.68%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '530.13'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +6.28%  '
$ws.Range('E40').Value = '  +13.85%  '
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('E42').Value = '  -0.02%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.861'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +8.86%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '3.70'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -4.80%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0421'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +20.58%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.31'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +4.92%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '5.49'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +15.04%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '8.29'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +12.71%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  +11.66%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.07'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +5.47%  '
